$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 0.5537637576255825
$ws.Range("D2").Value = 0.5833643669488036

# Row 3
$ws.Range("C3").Value = -0.5086014908303759
$ws.Range("D3").Value = 0.6143161562620532

# Row 4
$ws.Range("C4").Value = -0.09764012866084447
$ws.Range("D4").Value = 0.9227913984545495

# Row 5
$ws.Range("C5").Value = 0.4906697927758693
$ws.Range("D5").Value = 0.6268120133200328

# Row 6
$ws.Range("C6").Value = -1.265061016659617
$ws.Range("D6").Value = 0.2144514792615118

# Row 7
$ws.Range("C7").Value = -0.5939888942621706
$ws.Range("D7").Value = 0.5564513505541528

# Row 8
$ws.Range("C8").Value = 0.07835147309814031
$ws.Range("D8").Value = 0.9380077867662306

# Row 9
$ws.Range("C9").Value = 0.369919842631242
$ws.Range("D9").Value = 0.7137365432727436

# Row 10
$ws.Range("C10").Value = 2.088287244325005
$ws.Range("D10").Value = 0.04433352515238353
$ws.Range("G10").Value = "Sí"

# Row 11
$ws.Range("C11").Value = 0.4248076588568984
$ws.Range("D11").Value = 0.6736532880805541
